$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 14.28763719528746
$ws.Cells.Item(2, 3).Value = 13.57739184992911
$ws.Cells.Item(2, 4).Value = 4.341156593273206
$ws.Cells.Item(2, 6).Value = 20.44944119857568
$ws.Cells.Item(2, 7).Value = 22.13268719808953
$ws.Cells.Item(2, 8).Value = 12.65824369638209
$ws.Cells.Item(2, 9).Value = 19.9821073716528
$ws.Cells.Item(2, 12).Value = 11.01622645986776
$ws.Cells.Item(2, 13).Value = 14.45939861326351
$ws.Cells.Item(2, 15).Value = 18.37409351165983

$ws.Cells.Item(3, 2).Value = 13.62222537998159
$ws.Cells.Item(3, 3).Value = 13.40173155081667
$ws.Cells.Item(3, 4).Value = 4.256085223386085
$ws.Cells.Item(3, 6).Value = 20.50259121829934
$ws.Cells.Item(3, 7).Value = 22.21854043359447
$ws.Cells.Item(3, 8).Value = 12.71372382683534
$ws.Cells.Item(3, 9).Value = 20.12479791732333
$ws.Cells.Item(3, 12).Value = 11.04037989549928
$ws.Cells.Item(3, 13).Value = 14.31882187748653
$ws.Cells.Item(3, 15).Value = 18.4655489495442

$ws.Cells.Item(4, 2).Value = 13.19629600588131
$ws.Cells.Item(4, 3).Value = 13.29311027979399
$ws.Cells.Item(4, 4).Value = 4.202453401185002
$ws.Cells.Item(4, 6).Value = 20.54230196721789
$ws.Cells.Item(4, 7).Value = 22.28177580724263
$ws.Cells.Item(4, 8).Value = 12.75025590998129
$ws.Cells.Item(4, 9).Value = 20.21724082664188
$ws.Cells.Item(4, 12).Value = 11.05703551618506
$ws.Cells.Item(4, 13).Value = 14.2332096074509
$ws.Cells.Item(4, 15).Value = 18.52683226330372

$ws.Cells.Item(5, 2).Value = 13.01855328849873
$ws.Cells.Item(5, 3).Value = 13.24869302301696
$ws.Cells.Item(5, 4).Value = 4.180264107764927
$ws.Cells.Item(5, 6).Value = 20.56025666370326
$ws.Cells.Item(5, 7).Value = 22.31017104378979
$ws.Cells.Item(5, 8).Value = 12.76576314134915
$ws.Cells.Item(5, 9).Value = 20.25612821101293
$ws.Cells.Item(5, 12).Value = 11.06428174025765
$ws.Cells.Item(5, 13).Value = 14.19852834517033
$ws.Cells.Item(5, 15).Value = 18.55309148884597

$ws.Cells.Item(6, 2).Value = 12.98879345204995
$ws.Cells.Item(6, 3).Value = 13.24130947156169
$ws.Cells.Item(6, 4).Value = 4.176559956110337
$ws.Cells.Item(6, 6).Value = 20.56334484131295
$ws.Cells.Item(6, 7).Value = 22.31504408180518
$ws.Cells.Item(6, 8).Value = 12.7683755514246
$ws.Cells.Item(6, 9).Value = 20.2626589296623
$ws.Cells.Item(6, 12).Value = 11.06551268530771
$ws.Cells.Item(6, 13).Value = 14.19278287973194
$ws.Cells.Item(6, 15).Value = 18.55752935261611

$ws.Cells.Item(7, 2).Value = 13.19391552427521
$ws.Cells.Item(7, 3).Value = 13.29251182138239
$ws.Cells.Item(7, 4).Value = 4.202155476708445
$ws.Cells.Item(7, 6).Value = 20.54253694485988
$ws.Cells.Item(7, 7).Value = 22.28214814808246
$ws.Cells.Item(7, 8).Value = 12.75046253544579
$ws.Cells.Item(7, 9).Value = 20.21776034931607
$ws.Cells.Item(7, 12).Value = 11.05713138325052
$ws.Cells.Item(7, 13).Value = 14.2327410088854
$ws.Cells.Item(7, 15).Value = 18.52718120369322

$ws.Cells.Item(8, 2).Value = 14.06191386547264
$ws.Cells.Item(8, 3).Value = 13.51700634724864
$ws.Cells.Item(8, 4).Value = 4.312124336539457
$ws.Cells.Item(8, 6).Value = 20.46629453462701
$ws.Cells.Item(8, 7).Value = 22.16009455039228
$ws.Cells.Item(8, 8).Value = 12.67686105552683
$ws.Cells.Item(8, 9).Value = 20.03030544521479
$ws.Cells.Item(8, 12).Value = 11.02417569272743
$ws.Cells.Item(8, 13).Value = 14.41079994719161
$ws.Cells.Item(8, 15).Value = 18.40456035125478

$ws.Cells.Item(9, 2).Value = 15.61922586867429
$ws.Cells.Item(9, 3).Value = 13.94947587515151
$ws.Cells.Item(9, 4).Value = 4.515957069728207
$ws.Cells.Item(9, 6).Value = 20.3732160367077
$ws.Cells.Item(9, 7).Value = 22.00501978156978
$ws.Cells.Item(9, 8).Value = 12.55211493541287
$ws.Cells.Item(9, 9).Value = 19.70096383719544
$ws.Cells.Item(9, 12).Value = 10.9740346252634
$ws.Cells.Item(9, 13).Value = 14.76422565474504
$ws.Cells.Item(9, 15).Value = 18.20497812000044

$ws.Cells.Item(10, 2).Value = 16.66737555131479
$ws.Cells.Item(10, 3).Value = 14.26026495983988
$ws.Cells.Item(10, 4).Value = 4.657603348809118
$ws.Cells.Item(10, 6).Value = 20.33956724233782
$ws.Cells.Item(10, 7).Value = 21.94342975622351
$ws.Cells.Item(10, 8).Value = 12.47241939760199
$ws.Cells.Item(10, 9).Value = 19.48222136971546
$ws.Cells.Item(10, 12).Value = 10.94602731333844
$ws.Cells.Item(10, 13).Value = 15.02470506283116
$ws.Cells.Item(10, 15).Value = 18.083505390073

$ws.Cells.Item(11, 2).Value = 17.12214873612405
$ws.Cells.Item(11, 3).Value = 14.39969184956747
$ws.Cells.Item(11, 4).Value = 4.720107983462225
$ws.Cells.Item(11, 6).Value = 20.33185403050544
$ws.Cells.Item(11, 7).Value = 21.9269462760385
$ws.Cells.Item(11, 8).Value = 12.43876345139876
$ws.Cells.Item(11, 9).Value = 19.38773328504634
$ws.Cells.Item(11, 12).Value = 10.93520296194277
$ws.Cells.Item(11, 13).Value = 15.14300559555516
$ws.Cells.Item(11, 15).Value = 18.03375823356702

$ws.Cells.Item(12, 2).Value = 17.29111228628535
$ws.Cells.Item(12, 3).Value = 14.45217215719448
$ws.Cells.Item(12, 4).Value = 4.743485970418673
$ws.Cells.Item(12, 6).Value = 20.33002785240195
$ws.Cells.Item(12, 7).Value = 21.92237388786886
$ws.Cells.Item(12, 8).Value = 12.42639279425754
$ws.Cells.Item(12, 9).Value = 19.35267361969347
$ws.Cells.Item(12, 12).Value = 10.93137946836037
$ws.Cells.Item(12, 13).Value = 15.18774399090876
$ws.Cells.Item(12, 15).Value = 18.01571694166497

$ws.Cells.Item(13, 2).Value = 17.25486860253583
$ws.Cells.Item(13, 3).Value = 14.4408842830862
$ws.Cells.Item(13, 4).Value = 4.738464256803073
$ws.Cells.Item(13, 6).Value = 20.33037243645856
$ws.Cells.Item(13, 7).Value = 21.92328422803436
$ws.Cells.Item(13, 8).Value = 12.42904038885303
$ws.Cells.Item(13, 9).Value = 19.36019229453261
$ws.Cells.Item(13, 12).Value = 10.93219067763209
$ws.Cells.Item(13, 13).Value = 15.17811192087215
$ws.Cells.Item(13, 15).Value = 18.01956694703053

$ws.Cells.Item(14, 2).Value = 17.13611499823887
$ws.Cells.Item(14, 3).Value = 14.40401600814167
$ws.Cells.Item(14, 4).Value = 4.722037213715777
$ws.Cells.Item(14, 6).Value = 20.33168183980437
$ws.Cells.Item(14, 7).Value = 21.92653658746967
$ws.Cells.Item(14, 8).Value = 12.43773820895625
$ws.Cells.Item(14, 9).Value = 19.38483446119933
$ws.Cells.Item(14, 12).Value = 10.93488288116795
$ws.Cells.Item(14, 13).Value = 15.14668762167709
$ws.Cells.Item(14, 15).Value = 18.03225797042227

$ws.Cells.Item(15, 2).Value = 17.06294970132668
$ws.Cells.Item(15, 3).Value = 14.3813906895866
$ws.Cells.Item(15, 4).Value = 4.711936873561618
$ws.Cells.Item(15, 6).Value = 20.33262650393701
$ws.Cells.Item(15, 7).Value = 21.92874646497413
$ws.Cells.Item(15, 8).Value = 12.44311461246869
$ws.Cells.Item(15, 9).Value = 19.40002236286798
$ws.Cells.Item(15, 12).Value = 10.93656780146016
$ws.Cells.Item(15, 13).Value = 15.12743065583869
$ws.Cells.Item(15, 15).Value = 18.04013549316789

$ws.Cells.Item(16, 2).Value = 16.63720442977077
$ws.Cells.Item(16, 3).Value = 14.25111061337708
$ws.Cells.Item(16, 4).Value = 4.653478492869232
$ws.Cells.Item(16, 6).Value = 20.34022437749523
$ws.Cells.Item(16, 7).Value = 21.94474007102464
$ws.Cells.Item(16, 8).Value = 12.47467121584062
$ws.Cells.Item(16, 9).Value = 19.4884973123267
$ws.Cells.Item(16, 12).Value = 10.94677326320515
$ws.Cells.Item(16, 13).Value = 15.01696724242625
$ws.Cells.Item(16, 15).Value = 18.08686776025795

$ws.Cells.Item(17, 2).Value = 16.37031771473173
$ws.Cells.Item(17, 3).Value = 14.17066126503659
$ws.Cells.Item(17, 4).Value = 4.617111700041328
$ws.Cells.Item(17, 6).Value = 20.34683260997928
$ws.Cells.Item(17, 7).Value = 21.95751425612613
$ws.Cells.Item(17, 8).Value = 12.49469594521841
$ws.Cells.Item(17, 9).Value = 19.54405870870291
$ws.Cells.Item(17, 12).Value = 10.95352473750833
$ws.Cells.Item(17, 13).Value = 14.94912941487648
$ws.Cells.Item(17, 15).Value = 18.11695120844931

$ws.Cells.Item(18, 2).Value = 16.214741269358
$ws.Cells.Item(18, 3).Value = 14.12420726158819
$ws.Cells.Item(18, 4).Value = 4.596013476324276
$ws.Cells.Item(18, 6).Value = 20.35134816010445
$ws.Cells.Item(18, 7).Value = 21.96594675164962
$ws.Cells.Item(18, 8).Value = 12.50645813854397
$ws.Cells.Item(18, 9).Value = 19.57648860507969
$ws.Cells.Item(18, 12).Value = 10.95758838985548
$ws.Cells.Item(18, 13).Value = 14.91009461094499
$ws.Cells.Item(18, 15).Value = 18.13477291569185

$ws.Cells.Item(19, 2).Value = 16.16171288228759
$ws.Cells.Item(19, 3).Value = 14.10844867692578
$ws.Cells.Item(19, 4).Value = 4.588839319012616
$ws.Cells.Item(19, 6).Value = 20.35299969454687
$ws.Cells.Item(19, 7).Value = 21.96898785241686
$ws.Cells.Item(19, 8).Value = 12.51048259187931
$ws.Cells.Item(19, 9).Value = 19.58754998041004
$ws.Cells.Item(19, 12).Value = 10.95899525607165
$ws.Cells.Item(19, 13).Value = 14.89687627471035
$ws.Cells.Item(19, 15).Value = 18.14089595943429

$ws.Cells.Item(20, 2).Value = 16.39894315502499
$ws.Cells.Item(20, 3).Value = 14.17924430832777
$ws.Cells.Item(20, 4).Value = 4.621001841445184
$ws.Cells.Item(20, 6).Value = 20.34605516406493
$ws.Cells.Item(20, 7).Value = 21.95604203020537
$ws.Cells.Item(20, 8).Value = 12.49253897005611
$ws.Cells.Item(20, 9).Value = 19.53809521238122
$ws.Cells.Item(20, 12).Value = 10.95278736320794
$ws.Cells.Item(20, 13).Value = 14.95635278352366
$ws.Cells.Item(20, 15).Value = 18.11369507914516

$ws.Cells.Item(21, 2).Value = 17.17108456356298
$ws.Cells.Item(21, 3).Value = 14.41485402160241
$ws.Cells.Item(21, 4).Value = 4.726870239297243
$ws.Cells.Item(21, 6).Value = 20.33126751180778
$ws.Cells.Item(21, 7).Value = 21.9255359028174
$ws.Cells.Item(21, 8).Value = 12.43517329051296
$ws.Cells.Item(21, 9).Value = 19.37757689984175
$ws.Cells.Item(21, 12).Value = 10.93408464117334
$ws.Cells.Item(21, 13).Value = 15.15591956390262
$ws.Cells.Item(21, 15).Value = 18.02850864797784

$ws.Cells.Item(22, 2).Value = 17.65675715254643
$ws.Cells.Item(22, 3).Value = 14.56697106844949
$ws.Cells.Item(22, 4).Value = 4.794358655836125
$ws.Cells.Item(22, 6).Value = 20.32798386281392
$ws.Cells.Item(22, 7).Value = 21.91533335642296
$ws.Cells.Item(22, 8).Value = 12.39986241563268
$ws.Cells.Item(22, 9).Value = 19.27687067822648
$ws.Cells.Item(22, 12).Value = 10.92346679153083
$ws.Cells.Item(22, 13).Value = 15.28598712841333
$ws.Cells.Item(22, 15).Value = 17.97748115630577

$ws.Cells.Item(23, 2).Value = 17.39930299721134
$ws.Cells.Item(23, 3).Value = 14.48596604017623
$ws.Cells.Item(23, 4).Value = 4.758498826377529
$ws.Cells.Item(23, 6).Value = 20.32915194266402
$ws.Cells.Item(23, 7).Value = 21.91988482949642
$ws.Cells.Item(23, 8).Value = 12.41850875925619
$ws.Cells.Item(23, 9).Value = 19.33023526251833
$ws.Cells.Item(23, 12).Value = 10.92898688958321
$ws.Cells.Item(23, 13).Value = 15.21661098910034
$ws.Cells.Item(23, 15).Value = 18.00428889480714

$ws.Cells.Item(24, 2).Value = 16.38600825197215
$ws.Cells.Item(24, 3).Value = 14.17536454270747
$ws.Cells.Item(24, 4).Value = 4.619243700756495
$ws.Cells.Item(24, 6).Value = 20.3464044160872
$ws.Cells.Item(24, 7).Value = 21.9567042339861
$ws.Cells.Item(24, 8).Value = 12.49351336073216
$ws.Cells.Item(24, 9).Value = 19.54078979230276
$ws.Cells.Item(24, 12).Value = 10.95312016277574
$ws.Cells.Item(24, 13).Value = 14.95308720182699
$ws.Cells.Item(24, 15).Value = 18.11516553603981

$ws.Cells.Item(25, 2).Value = 15.2143193939933
$ws.Cells.Item(25, 3).Value = 13.83354448168468
$ws.Cells.Item(25, 4).Value = 4.462177107619338
$ws.Cells.Item(25, 6).Value = 20.39231671605425
$ws.Cells.Item(25, 7).Value = 22.03784282689373
$ws.Cells.Item(25, 8).Value = 12.58376372669706
$ws.Cells.Item(25, 9).Value = 19.785974185584
$ws.Cells.Item(25, 12).Value = 10.98604810416971
$ws.Cells.Item(25, 13).Value = 14.66834263221202
$ws.Cells.Item(25, 15).Value = 18.25456920160135
